$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Top three" header in J1 and the SUM formulas in J2:J10
# (leftover leaderboard column no longer needed)
$ws.Range("J1:J10").ClearContents()

# Update the active selection to F12
$ws.Range("F12").Select()
